$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Automated refresh: shift the timestamp blocks in column D down one slot
# and stamp the newest refresh time into the first block (rows 2-15).
# Using .Formula with a plain numeric literal (no leading "=") keeps the
# cell a plain numeric value while preserving full double precision.
$ws.Range("D2:D15").Formula = "44232.53506424778"
$ws.Range("D16:D29").Formula = "44232.51400987268"
$ws.Range("D30:D37").Formula = "44232.49295149306"
